# aggiornamento fino a 6 gennaio 2022
# Append the new daily rows (465:491) to the Sassuolo report and extend
# the sheet dimension accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(465, 44539, 7, 47, 116.6890113709718)
    ,@(466, 44540, 13, 54, 134.0682258304782)
    ,@(467, 44541, 6, 49, 121.654501216545)
    ,@(468, 44542, 14, 55, 136.5509707532648)
    ,@(469, 44543, 11, 59, 146.4819504444114)
    ,@(470, 44544, 19, 70, 173.7921445950643)
    ,@(471, 44545, 2, 72, 178.7576344406376)
    ,@(472, 44546, 18, 83, 206.0678285912905)
    ,@(473, 44547, 22, 92, 228.4125328963702)
    ,@(474, 44548, 22, 108, 268.1364516609564)
    ,@(475, 44550, 32, 126, 312.8258602711157)
    ,@(476, 44551, 28, 143, 355.0325239584885)
    ,@(477, 44552, 9, 133, 330.2050747306222)
    ,@(478, 44553, 32, 163, 404.6874224142211)
    ,@(479, 44554, 33, 178, 441.9285962560206)
    ,@(480, 44555, 33, 189, 469.2387904066737)
    ,@(481, 44556, 52, 219, 543.7211380902726)
    ,@(482, 44557, 44, 231, 573.5140771637122)
    ,@(483, 44558, 13, 216, 536.2729033219127)
    ,@(484, 44559, 40, 247, 613.2379959282983)
    ,@(485, 44560, 86, 301, 747.3062217587766)
    ,@(486, 44561, 139, 407, 1010.47718357416)
    ,@(487, 44562, 103, 477, 1184.269328169224)
    ,@(488, 44563, 27, 452, 1122.200705099558)
    ,@(489, 44564, 23, 431, 1070.063061721039)
    ,@(490, 44565, 9, 427, 1060.132082029892)
    ,@(491, 44566, 31, 418, 1037.787377724813)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Column A holds dates; give the new cells the same style (date number
# format + border + alignment, cellXf index 2) as the rest of the column
# by copying the format from the last pre-existing date cell (A464).
$ws.Range("A464").Copy() | Out-Null
$ws.Range("A465:A491").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) only touches formatting, but make sure the
# date values themselves are still intact afterwards.
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
}
